# Apply the "FV2304/FV2310" column header rename + table + frozen header row
# change described by the commit:
#   "chore: adapt column header formatting to respective input file names"
#
# The worksheet has 21 header columns in row 1:
#   A..J  = "<Name>_old"  -> "<Name>_FV2304"
#   K     = "diff"        -> unchanged
#   L..U  = "<Name>_new"  -> "<Name>_FV2310"
# It is then wrapped in an Excel Table (ListObject) and the header row is
# frozen.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headers = @(
    "Segmentname_FV2304","Segmentgruppe_FV2304","Segment_FV2304","Datenelement_FV2304","Segment ID_FV2304",
    "Code_FV2304","Qualifier_FV2304","Beschreibung_FV2304","Bedingungsausdruck_FV2304","Bedingung_FV2304",
    "diff",
    "Segmentname_FV2310","Segmentgruppe_FV2310","Segment_FV2310","Datenelement_FV2310","Segment ID_FV2310",
    "Code_FV2310","Qualifier_FV2310","Beschreibung_FV2310","Bedingungsausdruck_FV2310","Bedingung_FV2310"
)

# 1) Rename the header cells (A1:U1) in place.
for ($i = 0; $i -lt $headers.Count; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# 2) Turn A1:U62 into an Excel Table named "Table1".
#    The header row already carries explicit (bold/filled/bordered)
#    formatting from the original file. If that formatting is still present
#    the instant the table is created, Excel captures it as a "manual
#    header override" (a dxf + headerRowDxfId on the table). To keep the
#    original header formatting AND avoid that extra dxf, we temporarily
#    stash the header formatting, reset the header range to the default
#    style, create the table, and then paste the stashed formatting back.
$headerRange = $ws.Range("A1:U1")
$stashRange  = $ws.Range("A100:U100")

$headerRange.Copy()
$stashRange.PasteSpecial(-4122)   # xlPasteFormats
$headerRange.Style = "Normal"

$rng = $ws.Range("A1:U62")
$tbl = $ws.ListObjects.Add(1, $rng, [System.Reflection.Missing]::Value, 1)
$tbl.Name = "Table1"
$tbl.TableStyle = ""

$stashRange.Copy()
$headerRange.PasteSpecial(-4122)  # xlPasteFormats
$stashRange.Clear()

# 3) Freeze the header row (split after row 1).
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
